$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row cell values (multiple primary keys handling)
$ws.Range("B1").Value = "Risk"
$ws.Range("C1").Value = "Cell"
$ws.Range("D1").Value = "Mask"

# Update the active selection on the sheet
$ws.Range("J11").Select()
